# Upload Leave Card 12/27/2023 4:01 PM
# Insert a new leave-card entry row (row 402) into Table1 on Sheet1, shifting
# all subsequent rows down by one, and update the new row plus the two rows
# immediately below it with the new SL (sick leave) entries.

$wb  = $excel.ActiveWorkbook
$ws  = $wb.Worksheets.Item("Sheet1")
$tbl = $ws.ListObjects.Item("Table1")

# --- 1. Insert a new row at 402, shifting existing rows (402..446) down to (403..447) ---
# Copy row 402 first so the newly inserted row inherits the same formatting
# (number formats, borders, table-formula) as the row that used to be there.
$ws.Rows.Item(402).Copy()
$ws.Rows.Item(402).Insert()

# Re-apply the original row-402 formatting (formats only) onto the new row 402,
# since Insert() on this engine does not reliably carry borders across.
$ws.Range("A403:K403").Copy()
$ws.Range("A402:K402").PasteSpecial(-4122)

# Grow the table definition so it covers the new last row (was K446, now K447).
$tbl.Resize($ws.Range("A8:K447"))

# --- 2. Populate the new row 402 (SL(1-0-0), 1 day, dated 9/29/2023) ---
$ws.Range("A402").Value = ""
$ws.Range("B402").Value = "SL(1-0-0)"
$ws.Range("C402").Value = ""
$ws.Range("H402").Value = 1
$ws.Range("K401").Copy()
$ws.Range("K402").PasteSpecial(-4122)
$ws.Range("K402").Value = 45198

# --- 3. Update row 403 (was old row 402) ---
$ws.Range("B403").Value = "SL(1-0-0)"
$ws.Range("C403").Value = 1.25
$ws.Range("H403").Value = 1
$ws.Range("K401").Copy()
$ws.Range("K403").PasteSpecial(-4122)
$ws.Range("K403").Value = 45230

# --- 4. Update row 404 (was old row 403) ---
$ws.Range("B404").Value = "SL(2-0-0)"
$ws.Range("H404").Value = 2
$ws.Range("K404").Value = "11/13,14/2023"

# --- 5. Re-enter the calculated-column formulas for the touched rows so the
#         table's calculated column recomputes against the resized range. ---
$formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'
$ws.Range("G402").Formula = $formula
$ws.Range("G403").Formula = $formula
$ws.Range("G404").Formula = $formula
$ws.Range("G447").Formula = $formula

$wb.Application.CalculateFull()
